$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths ----
# ColumnWidth assignment via COM re-derives the stored OOXML width with a
# constant +5/6 char offset baked in by the pixel-rounding path, so we
# subtract it here to land on the exact target widths from the diff.
$offset = 5 / 6
$widths = @(32, 16, 28, 27, 70, 28, 25, 21, 38, 24)
for ($c = 1; $c -le $widths.Length; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c - 1] - $offset
}

# ---- Header row (row 1) ----
# Extend the header formatting (bold / border / centered style used by
# A1:C1) across the new header cells before writing their text.
$ws.Range("C1").Copy()
$ws.Range("D1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$headers = @(
    "Docente",
    "Situação",
    "Nícolas Oliveira de Araújo",
    "Ítalo Moraes Rocha Guedes",
    "Toshik Iarley da Silva",
    "Genaina Aparecida de Souza",
    "Josimar Aleixo da Silva",
    "Ariana Mota Pereira",
    "Marilia Cecilia de Souza Bittencourt",
    "Vinícius Martins Silva"
)
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# ---- Data rows (row 2 onward) ----
# Each entry: Docente, Situação, then citation text keyed by column index
# (3..10) for whichever columns have a citation in that row.
$data = @(
    @{ Row=2;  A="Adriano Perin";                   B="1 citação(ões)"; Cols=@{ 10="PERIN, A. (1)" } },
    @{ Row=3;  A="Ana Maria Mapeli";                 B="Ok";             Cols=@{} },
    @{ Row=4;  A="André Samuel Strassburger";        B="Ok";             Cols=@{} },
    @{ Row=5;  A="Arthur Bernardes Cecílio Filho";   B="4 citação(ões)"; Cols=@{ 3="FILHO, A (1)"; 7="FILHO, A (1)"; 8="FILHO, A (1)"; 9="FILHO, A (1)" } },
    @{ Row=6;  A="Diego Ismael Rocha";               B="Ok";             Cols=@{} },
    @{ Row=7;  A="Diego Silva Batista";              B="4 citação(ões)"; Cols=@{ 5="BATISTA, D. S. (1) | BATISTA, DIEGO S (1) | BATISTA, DIEGO SILVA (1)"; 6="BATISTA, D. S. (1)" } },
    @{ Row=8;  A="Fernando Cesar Sala";              B="Ok";             Cols=@{} },
    @{ Row=9;  A="Guilherme da Silva Pereira";       B="Ok";             Cols=@{} },
    @{ Row=10; A="Jackson Mirellys Azevêdo Souza";   B="1 citação(ões)"; Cols=@{ 4="SOUZA, J. M. A. (1)" } },
    @{ Row=11; A="Kassio Ferreira Mendes";           B="Ok";             Cols=@{} },
    @{ Row=12; A="Leilson Costa Grangeiro";          B="Ok";             Cols=@{} },
    @{ Row=13; A="Luis Felipe Villani Purquerio";    B="Ok";             Cols=@{} },
    @{ Row=14; A="Rumy Goto";                        B="Ok";             Cols=@{} },
    @{ Row=15; A="Simone da Costa Mello";            B="1 citação(ões)"; Cols=@{ 4="MELLO, S. C. (1)" } },
    @{ Row=16; A="Thiago de Oliveira Vargas";        B="Ok";             Cols=@{} }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    for ($c = 3; $c -le 10; $c++) {
        if ($entry.Cols.ContainsKey($c)) {
            $ws.Cells.Item($r, $c).Value = $entry.Cols[$c]
        } else {
            $ws.Cells.Item($r, $c).Value = ""
            $ws.Cells.Item($r, $c).NumberFormat = "@"
        }
    }
}
